# Add the "Banner" worksheet, positioned immediately before "WorkBasketResultFields"
# (i.e. right after "SearchResultFields"), matching the target workbook layout.
$wb = $excel.ActiveWorkbook

$refSheet = $wb.Worksheets.Item("WorkBasketResultFields")
$banner = $wb.Worksheets.Add($refSheet)
$banner.Name = "Banner"

# The engine seeds new sheets with a copy of another sheet's content/format;
# wipe that out so we start from a clean slate.
$banner.Cells.ClearContents()
$banner.Cells.ClearFormats()
$banner.StandardHeight = 16

# ---- Column widths ----
# (engine's ColumnWidth setter adds a fixed 5/6-character padding internally,
#  so the input is pre-compensated to land on the target stored width)
$pad = 0.8333333333333333
$banner.Columns.Item(1).ColumnWidth = 31 - $pad
$banner.Columns.Item(2).ColumnWidth = 31.83203125 - $pad
$banner.Columns.Item(3).ColumnWidth = 31.6640625 - $pad
$banner.Columns.Item(4).ColumnWidth = 30.5 - $pad

# ---- Row 1 : sheet title band ----
$banner.Range("A1").Value = "Banner"
$titleRow = $banner.Range("A1:D1")
$titleRow.Interior.Color = 13434879
$titleRow.Borders.LineStyle = 1
$banner.Range("A1").Font.Bold = $true
$banner.Range("A1").Font.Name = "Arial"
$banner.Range("A1").Font.Size = 10
$banner.Range("B1:D1").Font.Name = "Arial"
$banner.Range("B1:D1").Font.Size = 10

# ---- Row 2 : column descriptions ----
$banner.Range("A2").Value = "Yes or No, to enable or disable the banner.                          MaxLength: 30"
$banner.Range("B2").Value = "Content to display in the banner         MaxLength: 300"
$banner.Range("C2").Value = "Optional link text in the banner                   MaxLength:50"
$banner.Range("D2").Value = "The URL for the link text.                                               Max Length: 1000"
$descRow = $banner.Range("A2:D2")
$descRow.Font.Name = "Arial"
$descRow.Font.Size = 10
$descRow.Font.Italic = $true
$descRow.Interior.Color = 13434879
$descRow.Borders.LineStyle = 1
$descRow.WrapText = $true
$banner.Rows.Item(2).RowHeight = 43

# ---- Row 3 : column / field names ----
$banner.Range("A3").Value = "BannerEnabled"
$banner.Range("B3").Value = "BannerDescription"
$banner.Range("C3").Value = "BannerUrlText"
$banner.Range("D3").Value = "BannerUrl"
$nameRow = $banner.Range("A3:D3")
$nameRow.Font.Name = "Arial"
$nameRow.Font.Size = 10
$nameRow.Font.Bold = $true
$nameRow.Interior.Color = 13434879
$nameRow.Borders.LineStyle = 1

# ---- Row 4 : first (empty) data-entry row ----
$row4 = $banner.Range("A4:D4")
$row4.Font.Name = "Arial"
$row4.Font.Size = 10
$row4.Borders.LineStyle = 1

# ---- Rows 5-11 : additional (empty) data-entry rows ----
$fillRows = $banner.Range("A5:D11")
$fillRows.Font.Name = "Arial"
$fillRows.Font.Size = 10

# ---- Activate the Banner tab (matches tabSelected moving to this sheet) ----
$banner.Activate()
$banner.Range("A1:XFD1048576").Select()
